# Update gh-pages to output generated at 456a3b4
#
# The three earliest events (2024-07-06) drop off the "展览" (sheet1) and
# "全部类型" (sheet4) listings (they've passed), every remaining row shifts
# up by three, the running index in column A is renumbered, and the
# "想去人数" (want-to-go count) column F is reset to 0 for every remaining
# row on every sheet. Sheet2 ("演出") keeps its single row but also has its
# F column reset to 0. Sheet3 ("本地生活") only has a header row and is
# untouched.

$wb = $excel.ActiveWorkbook

function Update-EventSheet([object]$ws) {
    # Drop the three oldest events (rows 2-4); everything below shifts up.
    $ws.Range("2:4").EntireRow.Delete()

    $lastRow = $ws.UsedRange.Rows.Count

    # Renumber the running index in column A (1, 2, 3, ...) and zero out
    # the "want to go" counter in column F for every remaining data row.
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
        $ws.Cells.Item($r, 6).Value = 0
    }
}

# 展览 (Exhibitions)
Update-EventSheet $wb.Worksheets.Item("展览")

# 演出 (Performances) - only row 2, just reset the want-to-go count.
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value = 0

# 全部类型 (All types) - same shift/renumber/reset treatment.
Update-EventSheet $wb.Worksheets.Item("全部类型")
